$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.429.99"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").Value = "3.630.12"
$ws.Range("E3").Value = "  +5.57%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'201.59"
$ws.Range("E5").Value = "  +11.35%  "
$ws.Range("D6").Value = "'581.28"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("D7").Value = "3.621.60"
$ws.Range("E7").Value = "  +5.37%  "
$ws.Range("E8").Value = "  +4.98%  "
$ws.Range("D10").Value = "'0.687"
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("D11").Value = "'60.80"
$ws.Range("E11").Value = "  +19.80%  "
$ws.Range("E12").Value = "  +7.42%  "
$ws.Range("E13").Value = "  +15.66%  "
$ws.Range("D14").Value = "'10.19"
$ws.Range("E14").Value = "  +9.08%  "
$ws.Range("D15").Value = "4.201.63"
$ws.Range("E15").Value = "  +5.41%  "
$ws.Range("D16").Value = "3.630.83"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "'19.34"
$ws.Range("E17").Value = "  +9.91%  "
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "'12.51"
$ws.Range("E19").Value = "  +7.91%  "
$ws.Range("D20").Value = "68.285.78"
$ws.Range("E20").Value = "  +5.44%  "
$ws.Range("E21").Value = "  +5.94%  "
$ws.Range("D22").Value = "'406.61"
$ws.Range("E22").Value = "  +8.22%  "
$ws.Range("D23").Value = "'13.02"
$ws.Range("E23").Value = "  +23.35%  "
$ws.Range("D24").Value = "'4.28"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("D25").Value = "'85.84"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("E26").Value = "  +19.14%  "
$ws.Range("E27").Value = "  +5.99%  "
$ws.Range("D28").Value = "'12.67"
$ws.Range("E28").Value = "  +7.79%  "
$ws.Range("D29").Value = "'6.15"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("E30").Value = "  +11.80%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  +13.87%  "
$ws.Range("D32").Value = "'31.95"
$ws.Range("E32").Value = "  +6.94%  "
$ws.Range("D33").Value = "'684.61"
$ws.Range("E33").Value = "  +14.08%  "
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("E35").Value = "  +5.86%  "
$ws.Range("D36").Value = "'64.01"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").Value = "'42.05"
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("D38").Value = "'0.419"
$ws.Range("E38").Value = "  +8.85%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "0.0₃0773"
$ws.Range("E40").Value = "  +9.59%  "
$ws.Range("D41").Value = "'3.20"
$ws.Range("E41").Value = "  +19.46%  "
$ws.Range("D42").Value = "3.208.23"
$ws.Range("E42").Value = "  +10.42%  "
$ws.Range("D43").Value = "'0.136"
$ws.Range("E43").Value = "  +6.53%  "
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  +13.61%  "
$ws.Range("D46").Value = "'2.89"
$ws.Range("E46").Value = "  +31.28%  "
$ws.Range("D47").Value = "'2.86"
$ws.Range("E47").Value = "  +16.49%  "
$ws.Range("E48").Value = "  +7.77%  "
$ws.Range("E49").Value = "  +5.03%  "
$ws.Range("E50").Value = "  +9.14%  "
$ws.Range("D51").Value = "'3.10"
$ws.Range("E51").Value = "  -0.41%  "
